$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the query table (Ratio_Identifiers_Update)
$lo = $ws.ListObjects.Item(1)

# Add a new calculated column "KG" to the table
$newCol = $lo.ListColumns.Add()

# Set header text directly (ListColumn.Name setter isn't reliable here)
$ws.Cells.Item(1, 7).Value = "KG"

# Fill in the calculated column formula for every data row, using the
# structured table reference so it matches what Excel itself writes when a
# calculated column is created.
$formula = "=Ratio_Identifiers_Update[[#This Row],[PecFinLengths]]/Ratio_Identifiers_Update[[#This Row],[fishStdLength]]"
$rowCount = $lo.ListRows.Count
for ($i = 1; $i -le $rowCount; $i++) {
    $ws.Cells.Item($i + 1, 7).Formula = $formula
}

# Restore the selection to reflect the cell the user ended up on after
# inserting the column.
$ws.Range("G15").Select()
